$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numeric-looking price strings (e.g. "308.18").
# Force text storage so COM does not silently coerce them to numbers,
# then restore the original (unstyled) cell style once all values are set.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.874.80'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '1.816.55'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '308.18'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('D7').Value = '0.4613'
$ws.Range('E7').Value = '  -2.65%  '
$ws.Range('D8').Value = '0.3644'
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('E9').Value = '  -3.26%  '
$ws.Range('D10').Value = '0.8574'
$ws.Range('D12').Value = '0.07506'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '5.320'
$ws.Range('E13').Value = '  -2.44%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '6.501'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = '91.70'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.718.91'
$ws.Range('E16').Value = '  -8.22%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').Value = '0.000008576'
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '26.908.90'
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = '14.41'
$ws.Range('E21').Value = '  -2.65%  '
$ws.Range('D22').Value = '5.124'
$ws.Range('E22').Value = '  -3.70%  '
$ws.Range('E23').Value = '  -1.80%  '
$ws.Range('D24').Value = '2.049.53'
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('D25').Value = '151.58'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').Value = '1.840'
$ws.Range('E26').Value = '  -3.76%  '
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('D28').Value = '2.072'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('D29').Value = '5.083'
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('D30').Value = '115.18'
$ws.Range('E30').Value = '  -2.09%  '
$ws.Range('D31').Value = '0.08854'
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').Value = '2.957'
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('D33').Value = '4.409'
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('E34').Value = '  -4.36%  '
$ws.Range('D35').Value = '0.7125'
$ws.Range('E35').Value = '  -5.98%  '
$ws.Range('D36').Value = '1.072'
$ws.Range('E36').Value = '  -2.94%  '
$ws.Range('E37').Value = '  -2.13%  '
$ws.Range('D38').Value = '2.403'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('D40').Value = '2.917'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('D41').Value = '7.133'
$ws.Range('E41').Value = '  -2.85%  '
$ws.Range('D42').Value = '0.5132'
$ws.Range('E42').Value = '  -3.72%  '
$ws.Range('D43').Value = '0.1620'
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('D44').Value = '8.151'
$ws.Range('E44').Value = '  -4.28%  '
$ws.Range('D45').Value = '0.4785'
$ws.Range('E45').Value = '  -2.58%  '
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = '103.00'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').Value = '10.02'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('D49').Value = '0.06289'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = '1.614'
$ws.Range('E50').Value = '  -3.74%  '
$ws.Range('D51').Value = '63.81'
$ws.Range('E51').Value = '  -2.97%  '

$ws.Range("D2:D51").Style = "Normal"
